# CR_Matthew_Flenet.xlsx - finalisation des mouvements de base du vaisseau
# et debut du codage du mouvement des asteroides.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New column F width (a new narrower column appears next to "duree")
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 12.4

# ---------------------------------------------------------------------
# Row 8 : fix the typo in the task description
# ---------------------------------------------------------------------
$ws.Range("C8").Value = "codage des mouvements du vaisseau avec son importation"

# ---------------------------------------------------------------------
# Row 10 : the placeholder "..." becomes a real duration
# ---------------------------------------------------------------------
$ws.Range("D10").Value = "30 minutes"

# ---------------------------------------------------------------------
# Row 11 : new completed task entry (vaisseau movements)
# ---------------------------------------------------------------------
$ws.Range("B11").Value = 43784
$ws.Range("C11").Value = "codage des mouvements du vaisseau"
$ws.Range("D11").Value = "1 heure 30 minutes"
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 12 : new completed task entry (asteroides parameters)
# ---------------------------------------------------------------------
$ws.Range("B12").Value = 43784
$ws.Range("C12").Value = "codage des parametres aléatoires ou non des astéroides"
$ws.Range("D12").Value = "1 heure"
$ws.Range("E10").Copy()
$ws.Range("E12").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 13 : new completed task entry (vaisseau movements, new week)
# ---------------------------------------------------------------------
$ws.Range("B13").Value = 43798
$ws.Range("C13").Value = "codage des mouvements du vaisseau"
$ws.Range("D13").Value = "1 heure"
$ws.Range("E4").Copy()
$ws.Range("E13").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 14 : new completed task entry (asteroides movements)
# ---------------------------------------------------------------------
$ws.Range("B14").Value = 43798
$ws.Range("C14").Value = "codage des mouvements des astéroides, rectiligne,vertical"
$ws.Range("D14").Value = "2 heure"
$ws.Range("E4").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# apply the same "week" date style used by the sibling rows (B7/B8/B9/B13)
$ws.Range("B9").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = 43798

# ---------------------------------------------------------------------
# Row 15 : only the date gets filled in, rest left blank
# ---------------------------------------------------------------------
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 43798

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Update the active selection to reflect where the author last worked
# ---------------------------------------------------------------------
$ws.Range("D14").Select()
